# Add the new "no_internet" column (column E) with header + data,
# and update the saved view state (window position, selection, scroll).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column header ---
$ws.Range("E1").Value = "no_internet"

# --- New column data (rows 2-77) ---
$noInternet = @(
    0, 0.03, 0.04, 0.04, 0.1, 0.13, 0.23, 0.1, 0.12, 0.12,
    0.18, 0.07, 0.15, 0.14, 0.02, 0.04, 0.13, 0.24, 0.23, 0.18,
    0.04, 0.17, 0.34, 0.32, 0.11, 0.11, 0.31, 0.3, 0.14, 0.04,
    0.37, 0.31, 0.23, 0.18, 0.12, 0.21, 0.12, 0.15, 0.15, 0.09,
    0.07, 0.16, 0.19, 0.27, 0.05, 0.18, 0.2, 0.14, 0.18, 0.19,
    0.2, 0.46, 0.38, 0.26, 0.07, 0.34, 0.21, 0.33, 0.16, 0.38,
    0.29, 0.09, 0.17, 0.12, 0.12, 0.24, 0.14, 0, 0.07, 0.2,
    0.2, 0.22, 0.13, 0.1, 0.31, 0.03
)

for ($i = 0; $i -lt $noInternet.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $noInternet[$i]
}

# --- Update view state: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E78").Select()

# --- Update workbook window position ---
$excel.ActiveWindow.Left = 15680
